# Update the HIGH/LOW/CLOSE/LTP/VOL/9:25-CLOSE figures in the "fo high low"
# sheet (Sheet1) for every stock row (rows 2-17), columns B through G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @{
    2  = @(857.2, 844.65, 848.05, 848.75, 44, 846.55)
    3  = @(908, 895.3, 900.75, 901.95, 27, 903.6)
    4  = @(46375.9, 46114.55, 46255, 46288.05, 17, 46175)
    5  = @(369.4, 362.4, 365.15, 365.55, 177, 364.65)
    6  = @(537.75, 529.1, 532.1, 533.3, 63, 534.9)
    7  = @(504.25, 496.25, 497.7, 497.9, 104, 503.65)
    8  = @(1000.75, 991.5, 995.2, 995.65, 172, 994)
    9  = @(719.4, 709.35, 712.95, 714.25, 41, 717.3)
    10 = @(20270, 20191, 20240.1, 20249.6, 38, 20225.1)
    11 = @(2472.75, 2459.25, 2460.9, 2465.7, 75, 2467.15)
    12 = @(602.7, 597.1, 599.25, 599.65, 234, 598.75)
    13 = @(886.3, 876.6, 878.75, 879.4, 24, 885.5)
    14 = @(638.5, 629.75, 636.5, 636.55, 249, 632.7)
    15 = @(133, 130.35, 132.85, 132.55, 899, 132.6)
    16 = @(3618, 3582.25, 3605.15, 3611.45, 22, 3599)
    17 = @(3271.6, 3240, 3260.05, 3260.35, 9, 3252)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
